# Reviewer requested the first column on Sheet1 be widened so the long
# food / item names are fully visible instead of being truncated.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Excel stores column width internally in "characters" of the workbook's
# default font (Aptos Narrow 11 here) and then re-derives the serialized
# pixel-based width from that. Setting the width in characters is how a
# user would do this interactively (drag the column border / Format >
# Column Width), so we set ColumnWidth on column A (the only column that
# changed) to produce the widened column seen in the saved workbook.
$ws.Columns.Item(1).ColumnWidth = 30.7
